$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture the existing "TAREA" explanation block (currently on rows 5-7)
# before it gets shifted down by one row to make room for the new note.
$tareaLabel   = $ws.Range("F5").Value2
$tareaText    = $ws.Range("G5").Value2
$pythonText   = $ws.Range("G6").Value2
$notebookText = $ws.Range("G7").Value2

# Remove the old F5/G5 cells entirely (they will be re-created one row below).
$ws.Range("F5:G5").Clear()

# Re-create the "TAREA:" label (bold) and its text one row down, on row 6.
$ws.Range("F6").Value = $tareaLabel
$ws.Range("F6").Font.Bold = $true

$ws.Range("G6").Value = $tareaText
$ws.Range("G7").Value = $pythonText
$ws.Range("G8").Value = $notebookText

# Add the new note about Cap periodicity on row 4.
$ws.Range("G4").Value = "Los Caps tienen periodicidad trimestral."

# Update the selected cell to match the final state.
$ws.Range("I16").Select()
